# Update cryptos list values (price + volume%) per Mon Mar 20 18:31:39 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text formatting before writing so that numeric-looking
# strings like "1.001" are stored as text (matching the original inlineStr cells)
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.959.79"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.764.91"
$ws.Range("E3").Value = "  -2.81%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "339.97"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").Value = "0.3770"
$ws.Range("E7").Value = "  -3.97%  "
$ws.Range("D8").Value = "0.3364"
$ws.Range("E8").Value = "  -3.66%  "
$ws.Range("D9").Value = "45.77"
$ws.Range("E9").Value = "  -5.34%  "
$ws.Range("D10").Value = "1.140"
$ws.Range("E10").Value = "  -5.19%  "
$ws.Range("D11").Value = "0.07176"
$ws.Range("E11").Value = "  -5.56%  "
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").Value = "22.61"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").Value = "6.219"
$ws.Range("E14").Value = "  -4.92%  "
$ws.Range("D15").Value = "7.207"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "1.760.39"
$ws.Range("E16").Value = "  -2.85%  "
$ws.Range("D17").Value = "0.00001054"
$ws.Range("E17").Value = "  -5.04%  "
$ws.Range("D18").Value = "0.06584"
$ws.Range("E18").Value = "  -2.02%  "
$ws.Range("D19").Value = "80.64"
$ws.Range("E19").Value = "  -5.31%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").Value = "17.00"
$ws.Range("E21").Value = "  -4.89%  "
$ws.Range("D22").Value = "6.288"
$ws.Range("E22").Value = "  -4.43%  "
$ws.Range("D23").Value = "27.868.68"
$ws.Range("D24").Value = "11.78"
$ws.Range("E24").Value = "  -8.41%  "
$ws.Range("D25").Value = "2.371"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("D26").Value = "152.74"
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "2.358"
$ws.Range("E27").Value = "  -8.28%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "19.77"
$ws.Range("E28").Value = "  -7.73%  "
$ws.Range("D29").Value = "1.293"
$ws.Range("E29").Value = "  -13.71%  "
$ws.Range("D30").Value = "1.954.33"
$ws.Range("E30").Value = "  -3.05%  "
$ws.Range("D31").Value = "131.51"
$ws.Range("E31").Value = "  -2.98%  "
$ws.Range("D32").Value = "4.026"
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").Value = "5.854"
$ws.Range("E33").Value = "  -4.88%  "
$ws.Range("D34").Value = "0.08727"
$ws.Range("D35").Value = "12.36"
$ws.Range("E35").Value = "  -7.34%  "
$ws.Range("D36").Value = "0.02351"
$ws.Range("E36").Value = "  -2.97%  "
$ws.Range("D37").Value = "0.6584"
$ws.Range("E37").Value = "  -5.40%  "
$ws.Range("D38").Value = "0.06229"
$ws.Range("E38").Value = "  -4.93%  "
$ws.Range("D39").Value = "5.164"
$ws.Range("E39").Value = "  -6.81%  "
$ws.Range("D40").Value = "0.2108"
$ws.Range("E40").Value = "  -5.73%  "
$ws.Range("E41").Value = "  -4.12%  "
$ws.Range("D42").Value = "1.445"
$ws.Range("E42").Value = "  -10.43%  "
$ws.Range("D43").Value = "8.062"
$ws.Range("E43").Value = "  -5.49%  "
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("D45").Value = "13.70"
$ws.Range("E45").Value = "  -7.03%  "
$ws.Range("D46").Value = "3.838"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("D47").Value = "0.6046"
$ws.Range("D48").Value = "129.85"
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("D49").Value = "2.006"
$ws.Range("E49").Value = "  -7.60%  "
$ws.Range("D50").Value = "0.07263"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "1.178"
$ws.Range("E51").Value = "  +0.91%  "

# Restore the default (unstyled) cell style on column D now that the values are set,
# so the cells keep matching the original (unstyled) look.
$ws.Range("D2:D51").Style = "Normal"
